# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" (exhibitions) sheet and the "全部类型" (all types) sheet.
# Both sheets list the same events; "全部类型" simply has one extra
# row (a concert entry) inserted above, so the target rows are offset
# by one between the two sheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row => new value, for the "展览" sheet
$exhibitUpdates = @{
    2  = 21
    4  = 82
    7  = 2653
    9  = 241
    10 = 100
    11 = 9739
    13 = 241
    15 = 11662
    16 = 11936
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row => new value, for the "全部类型" sheet
$allUpdates = @{
    2  = 21
    4  = 82
    7  = 2653
    10 = 241
    11 = 100
    12 = 9739
    14 = 241
    16 = 11662
    17 = 11936
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
